# Add a hyperlinked GitHub-repo textbox ("Rectangle 4") to the codingTime() slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(43)   # "codingTime();" slide (last slide)

# The slide currently has shape ids 1 (group), 2 (Title) and 3 (Footer), so a
# freshly-added autoshape would normally come back as id=4 / "Rectangle 3".
# The target deck's shape is id=5 / "Rectangle 4", so burn one id/name by
# adding and immediately removing a throwaway shape first.
$throwaway = $s.Shapes.AddShape(1, 0, 0, 1, 1)
$throwaway.Delete()

# EMU -> points (1 pt = 12700 EMU):
#   off  x=2280860  y=5055087
#   ext  cx=5458546 cy=369332
$rect = $s.Shapes.AddShape(1, 179.59527559055118, 398.0383464566929, 429.8067716535433, 29.081259842519685)

$tf = $rect.TextFrame
$tr = $tf.TextRange
$tr.Text = "https://github.com/domenic/understanding-node"

# Hyperlink the whole run to the repo.
$action = $tr.ActionSettings(1)
$action.Hyperlink.Address = "https://github.com/domenic/understanding-node"

# Split into the three runs seen in the target markup: "https", "://", and
# the rest, by touching each sub-range individually.
$r1 = $tr.Characters(1, 5)
$r1.ActionSettings(1).Hyperlink.Address = "https://github.com/domenic/understanding-node"

$r2 = $tr.Characters(6, 3)
$r2.ActionSettings(1).Hyperlink.Address = "https://github.com/domenic/understanding-node"

$r3 = $tr.Characters(9, 38)
$r3.ActionSettings(1).Hyperlink.Address = "https://github.com/domenic/understanding-node"

# bodyPr wrap="none" + spAutoFit, matching the target shape's text box sizing.
$tf.WordWrap = 0
$tf.AutoSize = 1
